# Generate Report for Handback
# Adds a second localized file (32385396-5055-45ef-8731-9774dff4db0e.md) to the
# Overview / zh-cn / de-de tables, and marks the existing file
# (06605d30-a0a8-472b-8fa3-6f10e3bae56c.md) as handed back too.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b39017e925f9d4d46f81adf0766e66157070c3a2/e2e/"
$oldFile = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$newFile = "32385396-5055-45ef-8731-9774dff4db0e.md"
$oldFileUrl = $repoBase + $oldFile
$newFileUrl = $repoBase + $newFile

$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: add row 3 for the new file
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2017-02-09 09:37:35"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", "e2e\" + $newFile) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: update row 2 (handback info for the old file) + add row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

# Row 2 -- the existing file has now also been handed back.
$wsZh.Range("C2").Value = $status
$wsZh.Range("L2").Value = "2017-02-09 09:39:38"
$wsZh.Range("L2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.230e561e228ff076c6c011518a567862c8c7c783.zh-cn.xlf"
$wsZh.Range("J2").Value = $oldFile
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $oldFileUrl, "", "", $oldFile) | Out-Null

# Row 3 -- brand-new file.
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newFile
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.zh-cn.xlf"
$wsZh.Range("H3").Value = "2017-02-09 09:37:11"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = $newFile
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsZh.Range("K3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.zh-cn.xlf"
$wsZh.Range("L3").Value = "2017-02-09 09:39:38"
$wsZh.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = ""
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "True"
$wsZh.Range("P3").Value = ""
$wsZh.Range("Q3").Value = "False"
$wsZh.Range("R3").Value = ""

# ---------------------------------------------------------------------------
# de-de sheet: update row 2 (handback info for the old file) + add row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

# Row 2 -- the existing file has now also been handed back.
$wsDe.Range("C2").Value = $status
$wsDe.Range("L2").Value = "2017-02-09 09:40:07"
$wsDe.Range("L2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.230e561e228ff076c6c011518a567862c8c7c783.de-de.xlf"
$wsDe.Range("J2").Value = $oldFile
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $oldFileUrl, "", "", $oldFile) | Out-Null

# Row 3 -- brand-new file.
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newFile
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.de-de.xlf"
$wsDe.Range("H3").Value = "2017-02-09 09:37:35"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = $newFile
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsDe.Range("K3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.de-de.xlf"
$wsDe.Range("L3").Value = "2017-02-09 09:40:07"
$wsDe.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = ""
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "True"
$wsDe.Range("P3").Value = ""
$wsDe.Range("Q3").Value = "False"
$wsDe.Range("R3").Value = ""

# ---------------------------------------------------------------------------
# Column widths widened for the longer status / filename text (per sheet)
# ---------------------------------------------------------------------------
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.9777050018311
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.9777050018311

$wsZh.Range("C1").EntireColumn.ColumnWidth = 29.9777050018311
$wsZh.Range("J1").EntireColumn.ColumnWidth = 40
$wsZh.Range("K1").EntireColumn.ColumnWidth = 40

$wsDe.Range("C1").EntireColumn.ColumnWidth = 29.9777050018311
$wsDe.Range("J1").EntireColumn.ColumnWidth = 40
$wsDe.Range("K1").EntireColumn.ColumnWidth = 40
